$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "iaest-measure:siglas"
$ws.Range("B4").Value = "medida"
$ws.Range("B5").Value = "xsd:string"

$ws.Range("B1").Font.Name = $ws.Range("A1").Font.Name
$ws.Range("B2").Font.Name = $ws.Range("A1").Font.Name
$ws.Range("B3").Font.Name = $ws.Range("A1").Font.Name
$ws.Range("B4").Font.Name = $ws.Range("A1").Font.Name
$ws.Range("B5").Font.Name = $ws.Range("A1").Font.Name

$ws.Rows.Item(6).Delete()
